# Auto-generated Excel COM-interop edit script
# Applies the Kujata_Profits sheet value corrections per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 1733.6
$ws.Cells.Item(125, 9).Value = 1616
$ws.Cells.Item(125, 11).Value = 14544
$ws.Cells.Item(125, 13).Value = -12084

$ws.Cells.Item(132, 8).Value = 6806930
$ws.Cells.Item(132, 9).Value = 7250605.5
$ws.Cells.Item(132, 11).Value = 21751816.5
$ws.Cells.Item(132, 13).Value = -21749286.5

$ws.Cells.Item(138, 8).Value = 479436.88
$ws.Cells.Item(138, 9).Value = 1502.3043
$ws.Cells.Item(138, 10).Value = 641091.25
$ws.Cells.Item(138, 11).Value = 4506.9129
$ws.Cells.Item(138, 12).Value = 1923273.75
$ws.Cells.Item(138, 13).Value = 633.0870999999997
$ws.Cells.Item(138, 14).Value = -1933553.75


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2351.25
$ws.Cells.Item(32, 9).Value = 2119.2576
$ws.Cells.Item(32, 10).Value = 10007
$ws.Cells.Item(32, 11).Value = 2119.2576
$ws.Cells.Item(32, 12).Value = 10007
$ws.Cells.Item(32, 13).Value = -1832.2576
$ws.Cells.Item(32, 14).Value = -10581

$ws.Cells.Item(74, 8).Value = 1257.9697
$ws.Cells.Item(74, 9).Value = 562.4545000000001
$ws.Cells.Item(74, 11).Value = 562.4545000000001
$ws.Cells.Item(74, 13).Value = 311.5454999999999

$ws.Cells.Item(77, 8).Value = 1257.9697
$ws.Cells.Item(77, 9).Value = 562.4545000000001
$ws.Cells.Item(77, 11).Value = 2812.2725
$ws.Cells.Item(77, 13).Value = 1555.7275

$ws.Cells.Item(132, 8).Value = 1964.15
$ws.Cells.Item(132, 9).Value = 1709.8276
$ws.Cells.Item(132, 11).Value = 5129.4828
$ws.Cells.Item(132, 13).Value = -2599.4828


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1608.3
$ws.Cells.Item(31, 9).Value = 1324.5333
$ws.Cells.Item(31, 10).Value = 2459.6
$ws.Cells.Item(31, 11).Value = 1324.5333
$ws.Cells.Item(31, 12).Value = 2459.6
$ws.Cells.Item(31, 13).Value = -1029.5333
$ws.Cells.Item(31, 14).Value = -3049.6

$ws.Cells.Item(34, 8).Value = 1608.3
$ws.Cells.Item(34, 9).Value = 1324.5333
$ws.Cells.Item(34, 10).Value = 2459.6
$ws.Cells.Item(34, 11).Value = 1324.5333
$ws.Cells.Item(34, 12).Value = 2459.6
$ws.Cells.Item(34, 13).Value = -1122.5333
$ws.Cells.Item(34, 14).Value = -2863.6

$ws.Cells.Item(58, 8).Value = 1417.1333
$ws.Cells.Item(58, 9).Value = 1194.2222
$ws.Cells.Item(58, 11).Value = 1194.2222
$ws.Cells.Item(58, 13).Value = -991.2221999999999

$ws.Cells.Item(122, 8).Value = 818.8946999999999
$ws.Cells.Item(122, 9).Value = 793.6429000000001
$ws.Cells.Item(122, 10).Value = 889.6
$ws.Cells.Item(122, 11).Value = 2380.9287
$ws.Cells.Item(122, 12).Value = 2668.8
$ws.Cells.Item(122, 13).Value = 69.07129999999961
$ws.Cells.Item(122, 14).Value = -7568.8

$ws.Cells.Item(132, 8).Value = 4222.3
$ws.Cells.Item(132, 9).Value = 5016.778
$ws.Cells.Item(132, 11).Value = 15050.334
$ws.Cells.Item(132, 13).Value = -12520.334

$ws.Cells.Item(134, 8).Value = 1968.069
$ws.Cells.Item(134, 9).Value = 2091
$ws.Cells.Item(134, 11).Value = 6273
$ws.Cells.Item(134, 13).Value = -3738

$ws.Cells.Item(136, 8).Value = 1417.1333
$ws.Cells.Item(136, 9).Value = 1194.2222
$ws.Cells.Item(136, 11).Value = 3582.6666
$ws.Cells.Item(136, 13).Value = -1032.6666


$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1215
$ws.Cells.Item(5, 9).Value = 1364.625
$ws.Cells.Item(5, 10).Value = 702
$ws.Cells.Item(5, 11).Value = 4093.875
$ws.Cells.Item(5, 12).Value = 2106
$ws.Cells.Item(5, 13).Value = -3981.875
$ws.Cells.Item(5, 14).Value = -2330

$ws.Cells.Item(14, 8).Value = 416.5
$ws.Cells.Item(14, 9).Value = 416.5
$ws.Cells.Item(14, 11).Value = 1249.5
$ws.Cells.Item(14, 13).Value = -1076.5

$ws.Cells.Item(34, 8).Value = 5001590.5
$ws.Cells.Item(34, 9).Value = 797
$ws.Cells.Item(34, 11).Value = 2391
$ws.Cells.Item(34, 13).Value = -2307

$ws.Cells.Item(44, 8).Value = 1670.7142
$ws.Cells.Item(44, 9).Value = 565
$ws.Cells.Item(44, 10).Value = 2500
$ws.Cells.Item(44, 11).Value = 1695
$ws.Cells.Item(44, 12).Value = 7500
$ws.Cells.Item(44, 13).Value = -1297
$ws.Cells.Item(44, 14).Value = -8296

$ws.Cells.Item(107, 8).Value = 502.9
$ws.Cells.Item(107, 9).Value = 299
$ws.Cells.Item(107, 10).Value = 525.55554
$ws.Cells.Item(107, 11).Value = 897
$ws.Cells.Item(107, 12).Value = 1576.66662
$ws.Cells.Item(107, 13).Value = 1023
$ws.Cells.Item(107, 14).Value = -5416.66662

$ws.Cells.Item(135, 8).Value = 1215
$ws.Cells.Item(135, 9).Value = 1364.625
$ws.Cells.Item(135, 10).Value = 702
$ws.Cells.Item(135, 11).Value = 12281.625
$ws.Cells.Item(135, 12).Value = 6318
$ws.Cells.Item(135, 13).Value = -9746.625
$ws.Cells.Item(135, 14).Value = -11388


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1256.6774
$ws.Cells.Item(122, 9).Value = 1356.4546
$ws.Cells.Item(122, 10).Value = 1012.7778
$ws.Cells.Item(122, 11).Value = 4069.3638
$ws.Cells.Item(122, 12).Value = 3038.3334
$ws.Cells.Item(122, 13).Value = -1619.3638
$ws.Cells.Item(122, 14).Value = -7938.3334

$ws.Cells.Item(126, 8).Value = 1763.6538
$ws.Cells.Item(126, 9).Value = 1484.7646
$ws.Cells.Item(126, 10).Value = 2290.4443
$ws.Cells.Item(126, 11).Value = 4454.293799999999
$ws.Cells.Item(126, 12).Value = 6871.3329
$ws.Cells.Item(126, 13).Value = -1984.293799999999
$ws.Cells.Item(126, 14).Value = -11811.3329


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 13).ClearContents()
$ws.Cells.Item(39, 14).ClearContents()

$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 13).ClearContents()
$ws.Cells.Item(41, 14).ClearContents()

$ws.Cells.Item(61, 8).Value = 1926.25
$ws.Cells.Item(61, 9).Value = 1833.3334
$ws.Cells.Item(61, 10).Value = 1982
$ws.Cells.Item(61, 11).Value = 1833.3334
$ws.Cells.Item(61, 12).Value = 1982
$ws.Cells.Item(61, 13).Value = -1631.3334
$ws.Cells.Item(61, 14).Value = -2386

$ws.Cells.Item(113, 8).Value = 1926.25
$ws.Cells.Item(113, 9).Value = 1833.3334
$ws.Cells.Item(113, 10).Value = 1982
$ws.Cells.Item(113, 11).Value = 1833.3334
$ws.Cells.Item(113, 12).Value = 1982
$ws.Cells.Item(113, 13).Value = 336.6666
$ws.Cells.Item(113, 14).Value = -6322

$ws.Cells.Item(122, 8).Value = 15632876
$ws.Cells.Item(122, 9).Value = 22737782
$ws.Cells.Item(122, 11).Value = 68213346
$ws.Cells.Item(122, 13).Value = -68210896


$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 444.83334
$ws.Cells.Item(107, 9).Value = 465.54544
$ws.Cells.Item(107, 11).Value = 1396.63632
$ws.Cells.Item(107, 13).Value = 523.3636799999999

$ws.Cells.Item(122, 8).Value = 17335384
$ws.Cells.Item(122, 9).Value = 20002058
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 60006174
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 13).Value = -60003724
$ws.Cells.Item(122, 14).Value = -10900

$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 14).ClearContents()

